$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 13 de Mayo de 2020 a las 00:05"

# Data refresh: country case statistics updated, plus Alemania/Brasil and Togo ranking changes
$ws.Range("B4").Value = 1406398
$ws.Range("C4").Value = 20564
$ws.Range("D4").Value = 280438
$ws.Range("E4").Value = 1042699
$ws.Range("F4").Value = 16472
$ws.Range("G4").Value = 1466
$ws.Range("H4").Value = 83261
$ws.Range("A10").Value = "Alemania"
$ws.Range("B10").Value = 173171
$ws.Range("C10").Value = 595
$ws.Range("D10").Value = 147200
$ws.Range("E10").Value = 18233
$ws.Range("F10").Value = 1539
$ws.Range("G10").Value = 77
$ws.Range("H10").Value = 7738
$ws.Range("A11").Value = "Brasil"
$ws.Range("B11").Value = 173141
$ws.Range("C11").Value = 3998
$ws.Range("D11").Value = 67384
$ws.Range("E11").Value = 93692
$ws.Range("F11").Value = 8318
$ws.Range("G11").Value = 440
$ws.Range("H11").Value = 12065
$ws.Range("B41").Value = 12272
$ws.Range("C41").Value = 659
$ws.Range("D41").Value = 2971
$ws.Range("E41").Value = 8808
$ws.Range("G41").Value = 14
$ws.Range("H41").Value = 493
$ws.Range("B87").Value = 1661
$ws.Range("C87").Value = 135
$ws.Range("D87").Value = 173
$ws.Range("E87").Value = 1408
$ws.Range("G87").Value = 6
$ws.Range("H87").Value = 80
$ws.Range("A146").Value = "Togo"
$ws.Range("B146").Value = 199
$ws.Range("C146").Value = 18
$ws.Range("D146").Value = 92
$ws.Range("E146").Value = 96
$ws.Range("F146").Value = 0
$ws.Range("H146").Value = 11
$ws.Range("A147").Value = "Martinica"
$ws.Range("D147").Value = 91
$ws.Range("E147").Value = 82
$ws.Range("F147").Value = 4
$ws.Range("H147").Value = 14
$ws.Range("A148").Value = "Islas Feroe"
$ws.Range("B148").Value = 187
$ws.Range("D148").Value = 187
$ws.Range("E148").Value = 0
$ws.Range("F148").Value = 0
$ws.Range("A149").Value = "Madagascar"
$ws.Range("B149").Value = 186
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 101
$ws.Range("E149").Value = 85
$ws.Range("F149").Value = 1
$ws.Range("H149").Value = 0
$ws.Range("A150").Value = "Suazilandia"
$ws.Range("B150").Value = 184
$ws.Range("C150").Value = 9
$ws.Range("D150").Value = 28
$ws.Range("E150").Value = 154
$ws.Range("H150").Value = 2
$ws.Range("B162").Value = 113
$ws.Range("C162").Value = 4
$ws.Range("E162").Value = 67
$ws.Range("B211").Value = 9
$ws.Range("C211").Value = 1
$ws.Range("E211").Value = 2
